$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Round existing recall/precision/F2 values in rows 2 and 3 to 3 decimals ---
$ws.Range("C2").Value = 0.277
$ws.Range("D2").Value = 0.523
$ws.Range("E2").Value = 0.305
$ws.Range("F2").Value = 0.438
$ws.Range("G2").Value = 0.64

$ws.Range("C3").Value = 0.723
$ws.Range("D3").Value = 0.316
$ws.Range("E3").Value = 0.575
$ws.Range("F3").Value = 0.489
$ws.Range("G3").Value = 0.303

# --- Add new row 4 ("Predicted" for drones in agriculture) ---
# Copy formatting from the existing "Predicted" row (row 2) so the new
# row's style matches (bold, bordered, centered/top-aligned) cells.
$ws.Range("A2").Copy($ws.Range("A4"))
$ws.Range("A4").Value = "Predicted"

$b4 = @"

"geospatial data" OR "aerial photography" OR "irrigation management" OR "soil analysis" OR "smart farming" OR "yield estimation" OR "crop monitoring" OR "agricultural innovation" OR "drone technology" OR "climate monitoring" OR "weed detection" OR "pesticide spraying" OR "land surveying" OR "agricultural robotics" OR "aerial imaging" OR "variable rate application" OR "field surveillance" OR "agricultural drone" OR "drone mapping" OR "drones in agriculture" OR "harvest prediction" OR "crop scouting" OR "livestock tracking" OR "crop health assessment" OR "farm management software"

"@
$ws.Range("B4").Value = $b4
# Setting a multi-line value triggers an automatic row-height estimate;
# AutoFit() re-measures (and drops the stale explicit height) so the row
# keeps its default (un-pinned) height, same as row 2's existing long text.
$ws.Rows.Item(4).AutoFit()

$ws.Range("C4").Value = 0.12
$ws.Range("D4").Value = 0.036
$ws.Range("E4").Value = 0.082
$ws.Range("F4").Value = 0.518
$ws.Range("G4").Value = 0.142

# --- Add new row 5 ("Baseline" for drones in agriculture) ---
$ws.Range("A3").Copy($ws.Range("A5"))
$ws.Range("A5").Value = "Baseline"

$ws.Range("B5").Value = '"Drones in Agriculture"'

$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0.592
$ws.Range("G5").Value = 0
